$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "B2" = 1.02
    "C2" = 1.039402199679229
    "D2" = 1.048932527732853
    "E2" = 1.0482127093648
    "F2" = 1.059753573177616
    "I2" = 1.043445844378464
    "J2" = 1.044494410417391
    "K2" = 1.051690994312071
    "L2" = 1.050973183424949
    "M2" = 1.062482217206266
    "N2" = 1.018758214742993
    "B3" = 1.02
    "C3" = 1.040206742504947
    "D3" = 1.049583733476009
    "E3" = 1.048929942165498
    "F3" = 1.06055624966144
    "I3" = 1.043651859734275
    "J3" = 1.044944869597427
    "K3" = 1.052154760891009
    "L3" = 1.051502661677676
    "M3" = 1.063099214376889
    "N3" = 1.018908282805671
    "B4" = 1.02
    "C4" = 1.040727960047371
    "D4" = 1.05000560981265
    "E4" = 1.049394983898513
    "F4" = 1.061076677073986
    "I4" = 1.043784115444525
    "J4" = 1.04523628502582
    "K4" = 1.052454664208267
    "L4" = 1.05184553943157
    "M4" = 1.063498836828653
    "N4" = 1.019005342429523
    "B5" = 1.02
    "C5" = 1.04094722761568
    "D5" = 1.050183085051479
    "E5" = 1.049590711517203
    "F5" = 1.061295712128273
    "I5" = 1.043839463496786
    "J5" = 1.045358779953443
    "K5" = 1.052580697906797
    "L5" = 1.051989748338518
    "M5" = 1.063666928302666
    "N5" = 1.019046135177377
    "B6" = 1.02
    "C6" = 1.040984052214867
    "D6" = 1.050212890813895
    "E6" = 1.049623588128476
    "F6" = 1.061332503522502
    "I6" = 1.043848741860892
    "J6" = 1.04537934640699
    "K6" = 1.052601856800018
    "L6" = 1.052013965304496
    "M6" = 1.063695156839241
    "N6" = 1.01905298379339
    "B7" = 1.02
    "C7" = 1.040730889333192
    "D7" = 1.05000798078255
    "E7" = 1.049397598342539
    "F7" = 1.061079602861747
    "I7" = 1.043784856000165
    "J7" = 1.045237921874039
    "K7" = 1.052456348457927
    "L7" = 1.051847466112016
    "M7" = 1.063501082522937
    "N7" = 1.019005887548502
    "B8" = 1.02
    "C8" = 1.039673968295996
    "D8" = 1.049152500529865
    "E8" = 1.04845490487891
    "F8" = 1.060024624427265
    "I8" = 1.043515685115172
    "J8" = 1.04464665720075
    "K8" = 1.051847763767161
    "L8" = 1.051152066277406
    "M8" = 1.06269065372753
    "N8" = 1.01880893986377
    "B9" = 1.02
    "C9" = 1.037816401328918
    "D9" = 1.047648969530024
    "E9" = 1.046801076143669
    "F9" = 1.058173691364575
    "I9" = 1.043033368641917
    "J9" = 1.043604358819591
    "K9" = 1.05077400933384
    "L9" = 1.049928817371241
    "M9" = 1.061265592789425
    "N9" = 1.018461573553314
    "B10" = 1.02
    "C10" = 1.036581400423749
    "D10" = 1.046649375164247
    "E10" = 1.045703560506613
    "F10" = 1.056945289736331
    "I10" = 1.042706495876434
    "J10" = 1.04290929377669
    "K10" = 1.050057351343763
    "L10" = 1.049114843924346
    "M10" = 1.060317685349159
    "N10" = 1.018229809782619
    "B11" = 1.02
    "C11" = 1.036047455068022
    "D11" = 1.046217218502966
    "E11" = 1.045229543551416
    "F11" = 1.056414722094148
    "I11" = 1.042563703549412
    "J11" = 1.042608292365149
    "K11" = 1.049746853852976
    "L11" = 1.048762765530204
    "M11" = 1.059907758913256
    "N11" = 1.018129415355205
    "B12" = 1.02
    "C12" = 1.035849248739988
    "D12" = 1.046056799545436
    "E12" = 1.045053657037831
    "F12" = 1.056217848882676
    "I12" = 1.042510476556397
    "J12" = 1.042496483330341
    "K12" = 1.049631495708061
    "L12" = 1.048632046204037
    "M12" = 1.059755574450345
    "N12" = 1.018092118992852
    "B13" = 1.02
    "C13" = 1.035891758995568
    "D13" = 1.046091205251791
    "E13" = 1.045091376911202
    "F13" = 1.056260069623407
    "I13" = 1.042521902400749
    "J13" = 1.042520466894021
    "K13" = 1.049656241554445
    "L13" = 1.048660083281662
    "M13" = 1.059788214864343
    "N13" = 1.018100119425848
    "B14" = 1.02
    "C14" = 1.036031068715656
    "D14" = 1.046203956103912
    "E14" = 1.045215000942167
    "F14" = 1.056398444334673
    "I14" = 1.042559307610774
    "J14" = 1.042599050263307
    "K14" = 1.049737318825492
    "L14" = 1.048751959032343
    "M14" = 1.059895177646273
    "N14" = 1.018126332533652
    "B15" = 1.02
    "C15" = 1.036116918617428
    "D15" = 1.046273439367859
    "E15" = 1.045291194294596
    "F15" = 1.056483728556308
    "I15" = 1.042582329367185
    "J15" = 1.042647467619626
    "K15" = 1.049787269866109
    "L15" = 1.048808574481837
    "M15" = 1.059961091662395
    "N15" = 1.018142482593205
    "B16" = 1.02
    "C16" = 1.036616854012718
    "D16" = 1.046678070352841
    "E16" = 1.045735045212248
    "F16" = 1.056980530145645
    "I16" = 1.042715946181596
    "J16" = 1.042929269654989
    "K16" = 1.050077954375751
    "L16" = 1.049138218293539
    "M16" = 1.060344902003348
    "N16" = 1.018236471845965
    "B17" = 1.02
    "C17" = 1.036930670619393
    "D17" = 1.046932066594221
    "E17" = 1.046013787761949
    "F17" = 1.057292520501867
    "I17" = 1.04279942529396
    "J17" = 1.043106028566923
    "K17" = 1.050260245974569
    "L17" = 1.049345097306817
    "M17" = 1.060585797568305
    "N17" = 1.018295418672736
    "B18" = 1.02
    "C18" = 1.037113793357118
    "D18" = 1.047080283232796
    "E18" = 1.046176490627478
    "F18" = 1.057474628186692
    "I18" = 1.042847996171034
    "J18" = 1.043209125711497
    "K18" = 1.050366556128487
    "L18" = 1.049465802652361
    "M18" = 1.060726358289756
    "N18" = 1.01832979754617
    "B19" = 1.02
    "C19" = 1.037176246792094
    "D19" = 1.047130832229483
    "E19" = 1.046231987892733
    "F19" = 1.057536744014089
    "I19" = 1.042864537021637
    "J19" = 1.043244278566262
    "K19" = 1.050402802138889
    "L19" = 1.049506966168627
    "M19" = 1.06077429436024
    "N19" = 1.018341519192006
    "B20" = 1.02
    "C20" = 1.036896992901478
    "D20" = 1.046904808471323
    "E20" = 1.045983869202559
    "F20" = 1.05725903353183
    "I20" = 1.042790481292483
    "J20" = 1.043087064344062
    "K20" = 1.050240689605707
    "L20" = 1.049322897378909
    "M20" = 1.06055994652659
    "N20" = 1.018289094627019
    "B21" = 1.02
    "C21" = 1.035990042016867
    "D21" = 1.046170750921745
    "E21" = 1.04517859163732
    "F21" = 1.056357690803509
    "I21" = 1.042548297879409
    "J21" = 1.042575909511236
    "K21" = 1.049713444283535
    "L21" = 1.048724902292822
    "M21" = 1.059863677537383
    "N21" = 1.018118613571975
    "B22" = 1.02
    "C22" = 1.035420528675928
    "D22" = 1.045709818253339
    "E22" = 1.044673350063754
    "F22" = 1.055792157774753
    "I22" = 1.042394942795963
    "J22" = 1.04225450539122
    "K22" = 1.049381797113909
    "L22" = 1.048349256923242
    "M22" = 1.059426372028915
    "N22" = 1.018011394352046
    "B23" = 1.02
    "C23" = 1.035722369361239
    "D23" = 1.045954110001404
    "E23" = 1.04494108623781
    "F23" = 1.056091845254702
    "I23" = 1.042476341757827
    "J23" = 1.042424889316767
    "K23" = 1.049557623007565
    "L23" = 1.048548361096021
    "M23" = 1.059658151145574
    "N23" = 1.018068236067112
    "B24" = 1.02
    "C24" = 1.036912210169252
    "D24" = 1.046917125043389
    "E24" = 1.045997387750313
    "F24" = 1.057274164453891
    "I24" = 1.042794523075975
    "J24" = 1.043095633470554
    "K24" = 1.050249526341112
    "L24" = 1.049332928458734
    "M24" = 1.060571627344056
    "N24" = 1.018291952202908
    "B25" = 1.02
    "C25" = 1.038296039561447
    "D25" = 1.048037190286374
    "E25" = 1.047227751284656
    "F25" = 1.058651232218391
    "I25" = 1.04315900202086
    "J25" = 1.043873858888208
    "K25" = 1.051051751735924
    "L25" = 1.050244794236454
    "M25" = 1.061633636818104
    "N25" = 1.018551410977593
}

foreach ($cell in $changes.Keys) {
    $ws.Range($cell).Value = $changes[$cell]
}